# Sync latest changes - UI refinements and reservation collector updates
#
# 1) Add the new "nullGoods" checkbox-label HTML snippet to cell E5
#    (wrapped, multi-line text -> new shared string + wrap-text cell style).
# 2) Row 5 grows taller to fit the wrapped 3-line text.
# 3) Selection / scroll position moves from F18 to E21 (no more F1 freeze-left anchor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E5: new "no goods" checkbox label markup -----------------------------
$nullGoodsLabel = "<label class=`"nullGoodsLabel`" style=`"margin:5px 3px; width: 100%;`">`n                                                <input type=`"checkbox`" class=`"nullGoods`"> 상품 없음`n                                            </label>"

$cellE5 = $ws.Range("E5")
$cellE5.Value = $nullGoodsLabel
$cellE5.WrapText = $true

# Row 5 needs to grow to show the 3 wrapped lines of the label text.
$ws.Rows(5).RowHeight = 49.5

# --- Move the live selection to E21 (was F18) ------------------------------
$ws.Range("E21").Select()
